$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 971
$ws.Range("I18").Value = 971
$ws.Range("K18").Value = 971
$ws.Range("M18").Value = -687
$ws.Range("H19").Value = 1659.4
$ws.Range("I19").Value = 1299.75
$ws.Range("K19").Value = 1299.75
$ws.Range("M19").Value = -1124.75
$ws.Range("H32").Value = 3500.5
$ws.Range("I32").Value = 1999
$ws.Range("K32").Value = 1999
$ws.Range("M32").Value = -1673
$ws.Range("H40").Value = 4436.375
$ws.Range("I40").Value = 4098.2
$ws.Range("K40").Value = 4098.2
$ws.Range("M40").Value = -3923.2
$ws.Range("H55").Value = 279.45456
$ws.Range("I55").Value = 246.14285
$ws.Range("K55").Value = 246.14285
$ws.Range("M55").Value = -32.14285000000001
$ws.Range("H88").Value = 336416.5
$ws.Range("I88").Value = 1001749.5
$ws.Range("J88").Value = 3750
$ws.Range("K88").Value = 1001749.5
$ws.Range("L88").Value = 3750
$ws.Range("M88").Value = -1001343.5
$ws.Range("N88").Value = -4562
$ws.Range("H91").Value = 336416.5
$ws.Range("I91").Value = 1001749.5
$ws.Range("J91").Value = 3750
$ws.Range("K91").Value = 1001749.5
$ws.Range("L91").Value = 3750
$ws.Range("M91").Value = -1000345.5
$ws.Range("N91").Value = -6558
$ws.Range("H113").Value = 7802.0835
$ws.Range("J113").Value = 12974.75
$ws.Range("L113").Value = 12974.75
$ws.Range("N113").Value = -19482.75
$ws.Range("H116").Value = 6500.2856
$ws.Range("I116").Value = 7000
$ws.Range("J116").Value = 6300.4
$ws.Range("K116").Value = 7000
$ws.Range("L116").Value = 6300.4
$ws.Range("M116").Value = -3558
$ws.Range("N116").Value = -13184.4
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H126").Value = 114999
$ws.Range("J126").Value = 114999
$ws.Range("L126").Value = 114999
$ws.Range("M126").Value = -124879
$ws.Range("H132").Value = 1372.75
$ws.Range("I132").Value = 1382.5128
$ws.Range("J132").Value = 992
$ws.Range("K132").Value = 4147.538399999999
$ws.Range("L132").Value = 2976
$ws.Range("M132").Value = -1617.538399999999
$ws.Range("N132").Value = -8036
$ws.Range("H133").Value = 67500
$ws.Range("J133").Value = 67500
$ws.Range("L133").Value = 67500
$ws.Range("N133").Value = -77620
$ws.Range("H141").Value = 5050.3125
$ws.Range("J141").Value = 8599.200000000001
$ws.Range("L141").Value = 25797.6
$ws.Range("N141").Value = -36157.60000000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 600
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 2631.4
$ws.Range("I32").Value = 2442.4614
$ws.Range("K32").Value = 2442.4614
$ws.Range("M32").Value = -2155.4614
$ws.Range("H110").Value = 2107.6
$ws.Range("I110").Value = 2107.6
$ws.Range("K110").Value = 2107.6
$ws.Range("M110").Value = -62.59999999999991
$ws.Range("H122").Value = 2161.75
$ws.Range("I122").Value = 2277.3125
$ws.Range("K122").Value = 6831.9375
$ws.Range("M122").Value = -4381.9375

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 600
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H99").Value = 2785.4614
$ws.Range("I99").Value = 2890.111
$ws.Range("J99").Value = 2550
$ws.Range("K99").Value = 2890.111
$ws.Range("L99").Value = 2550
$ws.Range("M99").Value = -1392.111
$ws.Range("N99").Value = -5546
$ws.Range("H105").Value = 4046.8
$ws.Range("I105").Value = 3897.8
$ws.Range("J105").Value = 4195.8
$ws.Range("K105").Value = 3897.8
$ws.Range("L105").Value = 4195.8
$ws.Range("M105").Value = -2150.8
$ws.Range("N105").Value = -7689.8
$ws.Range("H134").Value = 1417.75
$ws.Range("I134").Value = 1417.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4253.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1718.25
$ws.Range("N134").Value = ""

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 562
$ws.Range("I22").Value = 550.8
$ws.Range("K22").Value = 550.8
$ws.Range("M22").Value = -200.8
$ws.Range("H31").Value = 3665.5557
$ws.Range("I31").Value = 2209
$ws.Range("K31").Value = 2209
$ws.Range("M31").Value = -1914
$ws.Range("H34").Value = 3665.5557
$ws.Range("I34").Value = 2209
$ws.Range("K34").Value = 2209
$ws.Range("M34").Value = -2007
$ws.Range("H86").Value = 4998.3335
$ws.Range("I86").Value = 4997.5
$ws.Range("K86").Value = 4997.5
$ws.Range("M86").Value = -3874.5
$ws.Range("H89").Value = 4998.3335
$ws.Range("I89").Value = 4997.5
$ws.Range("K89").Value = 24987.5
$ws.Range("M89").Value = -19371.5
$ws.Range("H134").Value = 2098.8
$ws.Range("I134").Value = 2132.5
$ws.Range("K134").Value = 6397.5
$ws.Range("M134").Value = -3862.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""
$ws.Range("H33").Value = 441
$ws.Range("J33").Value = 183
$ws.Range("L33").Value = 1098
$ws.Range("N33").Value = -1664
$ws.Range("H68").Value = 737.1
$ws.Range("I68").Value = 712.8570999999999
$ws.Range("K68").Value = 2138.5713
$ws.Range("M68").Value = -1327.5713
$ws.Range("H71").Value = 737.1
$ws.Range("I71").Value = 712.8570999999999
$ws.Range("K71").Value = 6415.7139
$ws.Range("M71").Value = -2359.7139
$ws.Range("H97").Value = 856
$ws.Range("I97").Value = 912
$ws.Range("K97").Value = 2736
$ws.Range("M97").Value = -2240

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 16258
$ws.Range("I40").Value = 16258
$ws.Range("K40").Value = 16258
$ws.Range("M40").Value = -16107
$ws.Range("H43").Value = 3439.625
$ws.Range("I43").Value = 2503.4
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 2503.4
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -2352.4
$ws.Range("N43").Value = -5302
$ws.Range("H102").Value = 4354.75
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 4354.75
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 4354.75
$ws.Range("M102").Value = ""
$ws.Range("N102").Value = -7598.75
$ws.Range("H113").Value = 3931.5
$ws.Range("I113").Value = 2518.2
$ws.Range("K113").Value = 2518.2
$ws.Range("M113").Value = -348.1999999999998

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2404.9412
$ws.Range("I7").Value = 1233.5
$ws.Range("J7").Value = 2765.3845
$ws.Range("K7").Value = 1233.5
$ws.Range("L7").Value = 2765.3845
$ws.Range("M7").Value = -1121.5
$ws.Range("N7").Value = -2989.3845
$ws.Range("H40").Value = 3317.375
$ws.Range("J40").Value = 3995
$ws.Range("L40").Value = 3995
$ws.Range("N40").Value = -4267
$ws.Range("H42").Value = 43333.332
$ws.Range("J42").Value = 60000
$ws.Range("L42").Value = 60000
$ws.Range("N42").Value = -61126
$ws.Range("H49").Value = 43333.332
$ws.Range("J49").Value = 60000
$ws.Range("L49").Value = 60000
$ws.Range("N49").Value = -60294
$ws.Range("H68").Value = 2543.1765
$ws.Range("J68").Value = 2179
$ws.Range("L68").Value = 2179
$ws.Range("N68").Value = -3677
$ws.Range("H71").Value = 2543.1765
$ws.Range("J71").Value = 2179
$ws.Range("L71").Value = 10895
$ws.Range("N71").Value = -18383
$ws.Range("H82").Value = 1566.4
$ws.Range("I82").Value = 1763.1818
$ws.Range("K82").Value = 1763.1818
$ws.Range("M82").Value = -1402.1818
$ws.Range("H85").Value = 1566.4
$ws.Range("I85").Value = 1763.1818
$ws.Range("K85").Value = 1763.1818
$ws.Range("M85").Value = -515.1818000000001
$ws.Range("H122").Value = 3713.6
$ws.Range("I122").Value = 3682.1765
$ws.Range("J122").Value = 3780.375
$ws.Range("K122").Value = 11046.5295
$ws.Range("L122").Value = 11341.125
$ws.Range("M122").Value = -8596.529500000001
$ws.Range("N122").Value = -16241.125
$ws.Range("H126").Value = 2404.9412
$ws.Range("I126").Value = 1233.5
$ws.Range("J126").Value = 2765.3845
$ws.Range("K126").Value = 3700.5
$ws.Range("L126").Value = 8296.1535
$ws.Range("M126").Value = -1230.5
$ws.Range("N126").Value = -13236.1535
$ws.Range("H132").Value = 2093.1428
$ws.Range("I132").Value = 1853.5
$ws.Range("K132").Value = 5560.5
$ws.Range("M132").Value = -3030.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1828.1428
$ws.Range("I81").Value = 1828.1428
$ws.Range("K81").Value = 3656.2856
$ws.Range("M81").Value = -2595.2856
$ws.Range("H84").Value = 1828.1428
$ws.Range("I84").Value = 1828.1428
$ws.Range("K84").Value = 18281.428
$ws.Range("M84").Value = -12977.428
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
